# AF-611: tests for named ranges and area intersections are moved to
# temporary_excel_files. The only functional change left inside this
# workbook itself is that the stray "=vfbjak" (#NAME?) formula in C5 is
# removed, which flips D7's ISBLANK(C5) result from FALSE to TRUE, and the
# active selection moves from D7 onto the now-empty C5.
# (The x15ac:absPath hint in workbook.xml and the customXml/SharePoint
# metadata part numbering are save-time/tooling artifacts with no Excel
# object-model surface, so they are out of scope for this script.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the erroneous formula cell C5 (=vfbjak -> #NAME?)
$ws.Range("C5").ClearContents()

# Move the saved selection to C5, matching the post-edit sheetView
$ws.Range("C5").Select()
